# Update "想去人数" (want-to-go count) figures in both the "展览" sheet and
# the "全部类型" rollup sheet, which mirrors the same rows.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 50
    "F4"  = 543
    "F9"  = 4497
    "F10" = 4373
    "F13" = 140
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
